$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per the financial data correction
$ws.Range("D2").Value = 274264
$ws.Range("E2").Value = 14055
$ws.Range("F2").Value = 14055
$ws.Range("G2").Value = 16035
$ws.Range("H2").Value = 13370
$ws.Range("I2").Value = 13375
$ws.Range("J2").Value = -5
$ws.Range("K2").Value = 2144337
$ws.Range("L2").Value = 1921490
$ws.Range("M2").Value = 222847
$ws.Range("N2").Value = 222049
$ws.Range("O2").Value = 798
$ws.Range("P2").Value = 1000
$ws.Range("Q2").Value = 89464
$ws.Range("R2").Value = -73923
$ws.Range("S2").Value = -3550
$ws.Range("T2").Value = 741
$ws.Range("V2").Value = 3853
$ws.Range("W2").Value = 5.13
$ws.Range("X2").Value = 4.88
$ws.Range("Y2").Value = 6.49
$ws.Range("Z2").Value = 0.66
$ws.Range("AA2").Value = 862.25
$ws.Range("AB2").Value = 23265.02
$ws.Range("AC2").Value = 6687
$ws.Range("AD2").Value = 17.42
$ws.Range("AE2").Value = 117440
$ws.Range("AF2").Value = 0.99
$ws.Range("AG2").Value = 1800
$ws.Range("AH2").Value = 1.55
$ws.Range("AI2").Value = 25.45
$ws.Range("AJ2").Value = 200000000

$ws.Range("D3").Value = 277059
$ws.Range("E3").Value = 11485
$ws.Range("F3").Value = 11485
$ws.Range("G3").Value = 13715
$ws.Range("H3").Value = 12096
$ws.Range("I3").Value = 12112
$ws.Range("J3").Value = -17
$ws.Range("K3").Value = 2303629
$ws.Range("L3").Value = 2066412
$ws.Range("M3").Value = 237217
$ws.Range("N3").Value = 236458
$ws.Range("O3").Value = 759
$ws.Range("P3").Value = 1000
$ws.Range("Q3").Value = 52881
$ws.Range("R3").Value = -71766
$ws.Range("S3").Value = -526
$ws.Range("T3").Value = 676
$ws.Range("V3").Value = 9327
$ws.Range("W3").Value = 4.15
$ws.Range("X3").Value = 4.37
$ws.Range("Y3").Value = 5.28
$ws.Range("Z3").Value = 0.54
$ws.Range("AA3").Value = 871.1
$ws.Range("AB3").Value = 25190.12
$ws.Range("AC3").Value = 6056
$ws.Range("AD3").Value = 18.16
$ws.Range("AE3").Value = 128139
$ws.Range("AF3").Value = 0.86
$ws.Range("AG3").Value = 1800
$ws.Range("AH3").Value = 1.64
$ws.Range("AI3").Value = 27.48
$ws.Range("AJ3").Value = 200000000

$ws.Range("D4").Value = 304286
$ws.Range("E4").Value = 9865
$ws.Range("F4").Value = 9865
$ws.Range("G4").Value = 26075
$ws.Range("H4").Value = 21500
$ws.Range("I4").Value = 20543
$ws.Range("J4").Value = 957
$ws.Range("K4").Value = 2646538
$ws.Range("L4").Value = 2363454
$ws.Range("M4").Value = 283084
$ws.Range("N4").Value = 266442
$ws.Range("O4").Value = 16642
$ws.Range("P4").Value = 1000
$ws.Range("Q4").Value = 48344
$ws.Range("R4").Value = -29341
$ws.Range("S4").Value = 5529
$ws.Range("T4").Value = 8252
$ws.Range("V4").Value = 127878
$ws.Range("W4").Value = 3.24
$ws.Range("X4").Value = 7.07
$ws.Range("Y4").Value = 8.550000000000001
$ws.Range("Z4").Value = 0.87
$ws.Range("AA4").Value = 834.89
$ws.Range("AB4").Value = 30318.82
$ws.Range("AC4").Value = 10271
$ws.Range("AD4").Value = 10.95
$ws.Range("AE4").Value = 148374
$ws.Range("AF4").Value = 0.76
$ws.Range("AG4").Value = 1200
$ws.Range("AH4").Value = 1.07
$ws.Range("AI4").Value = 10.49
$ws.Range("AJ4").Value = 200000000

$ws.Range("D5").Value = 319590
$ws.Range("E5").Value = 16906
$ws.Range("F5").Value = 16906
$ws.Range("G5").Value = 16828
$ws.Range("H5").Value = 12632
$ws.Range("I5").Value = 11661
$ws.Range("J5").Value = 971
$ws.Range("K5").Value = 2827138
$ws.Range("L5").Value = 2515922
$ws.Range("M5").Value = 311216
$ws.Range("N5").Value = 293589
$ws.Range("O5").Value = 17626
$ws.Range("P5").Value = 1000
$ws.Range("Q5").Value = 16082
$ws.Range("R5").Value = -38353
$ws.Range("S5").Value = 15582
$ws.Range("T5").Value = 3388
$ws.Range("V5").Value = 138461
$ws.Range("W5").Value = 5.29
$ws.Range("X5").Value = 3.95
$ws.Range("Y5").Value = 4.51
$ws.Range("Z5").Value = 0.46
$ws.Range("AA5").Value = 808.42
$ws.Range("AB5").Value = 33131.94
$ws.Range("AC5").Value = 5831
$ws.Range("AD5").Value = 21.35
$ws.Range("AE5").Value = 163491
$ws.Range("AF5").Value = 0.76
$ws.Range("AG5").Value = 2000
$ws.Range("AH5").Value = 1.61
$ws.Range("AI5").Value = 30.8
$ws.Range("AJ5").Value = 200000000

$ws.Range("D6").Value = 322409
$ws.Range("E6").Value = 25833
$ws.Range("F6").Value = 25833
$ws.Range("G6").Value = 23657
$ws.Range("H6").Value = 17337
$ws.Range("I6").Value = 16644
$ws.Range("K6").Value = 2894277
$ws.Range("L6").Value = 2589222
$ws.Range("M6").Value = 305055
$ws.Range("N6").Value = 289021
$ws.Range("P6").Value = 1000
$ws.Range("Q6").Value = 19561
$ws.Range("R6").Value = -44110
$ws.Range("S6").Value = 18329
$ws.Range("T6").Value = 14434
$ws.Range("V6").Value = 156099
$ws.Range("W6").Value = 8.01
$ws.Range("X6").Value = 5.38
$ws.Range("Y6").Value = 5.95
$ws.Range("Z6").Value = 0.61
$ws.Range("AA6").Value = 848.77
$ws.Range("AB6").Value = 32515.83
$ws.Range("AC6").Value = 8322
$ws.Range("AD6").Value = 9.81
$ws.Range("AE6").Value = 160947
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 2650
$ws.Range("AH6").Value = 3.25
$ws.Range("AI6").Value = 28.59
$ws.Range("AJ6").Value = 200000000

$ws.Range("D7").Value = 165374
$ws.Range("E7").Value = -116
$ws.Range("G7").Value = 15281
$ws.Range("H7").Value = 10995
$ws.Range("I7").Value = 10387
$ws.Range("K7").Value = 3065367
$ws.Range("L7").Value = 2692462
$ws.Range("M7").Value = 372903
$ws.Range("N7").Value = 353576
$ws.Range("P7").Value = 1000
$ws.Range("W7").Value = -0.07000000000000001
$ws.Range("X7").Value = 6.65
$ws.Range("Y7").Value = 3.23
$ws.Range("Z7").Value = 0.37
$ws.Range("AA7").Value = 722.03
$ws.Range("AC7").Value = 5194
$ws.Range("AD7").Value = 14.11
$ws.Range("AE7").Value = 196896
$ws.Range("AF7").Value = 0.37
$ws.Range("AG7").Value = 2648
$ws.Range("AH7").Value = 3.61
$ws.Range("AI7").Value = 50.99

$ws.Range("D8").Value = 165471
$ws.Range("E8").Value = 1308
$ws.Range("G8").Value = 17332
$ws.Range("H8").Value = 12807
$ws.Range("I8").Value = 11690
$ws.Range("K8").Value = 3191007
$ws.Range("L8").Value = 2806248
$ws.Range("M8").Value = 384758
$ws.Range("N8").Value = 365214
$ws.Range("P8").Value = 1000
$ws.Range("W8").Value = 0.79
$ws.Range("X8").Value = 7.74
$ws.Range("Y8").Value = 3.25
$ws.Range("Z8").Value = 0.41
$ws.Range("AA8").Value = 729.35
$ws.Range("AC8").Value = 5845
$ws.Range("AD8").Value = 11.89
$ws.Range("AE8").Value = 203377
$ws.Range("AF8").Value = 0.34
$ws.Range("AG8").Value = 2722
$ws.Range("AH8").Value = 3.92
$ws.Range("AI8").Value = 46.56

$ws.Range("D9").Value = 165337
$ws.Range("E9").Value = 1820
$ws.Range("G9").Value = 18077
$ws.Range("H9").Value = 13350
$ws.Range("I9").Value = 12242
$ws.Range("K9").Value = 3296647
$ws.Range("L9").Value = 2899889
$ws.Range("M9").Value = 396756
$ws.Range("N9").Value = 368402
$ws.Range("P9").Value = 1000
$ws.Range("W9").Value = 1.1
$ws.Range("X9").Value = 8.07
$ws.Range("Y9").Value = 3.34
$ws.Range("Z9").Value = 0.41
$ws.Range("AA9").Value = 730.9
$ws.Range("AC9").Value = 6121
$ws.Range("AD9").Value = 11.35
$ws.Range("AE9").Value = 205153
$ws.Range("AF9").Value = 0.34
$ws.Range("AG9").Value = 3124
$ws.Range("AH9").Value = 4.49
$ws.Range("AI9").Value = 51.03

# Remove obsolete FCF / cashflow-adjacent cells no longer reported for these periods
$ws.Range("U2").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("U6").ClearContents()
$ws.Range("Q7:U7").ClearContents()
$ws.Range("Q8:U8").ClearContents()
$ws.Range("Q9:U9").ClearContents()
